$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 2932.25
$ws.Range("I86").Value = 2167.818
$ws.Range("J86").Value = 3866.5557
$ws.Range("K86").Value = 2167.818
$ws.Range("L86").Value = 3866.5557
$ws.Range("M86").Value = -1044.818
$ws.Range("N86").Value = -6112.5557
$ws.Range("H89").Value = 2932.25
$ws.Range("I89").Value = 2167.818
$ws.Range("J89").Value = 3866.5557
$ws.Range("K89").Value = 10839.09
$ws.Range("L89").Value = 19332.7785
$ws.Range("M89").Value = -5223.09
$ws.Range("N89").Value = -30564.7785
$ws.Range("H106").Value = 67031.414
$ws.Range("I106").Value = 2050
$ws.Range("J106").Value = 80027.7
$ws.Range("K106").Value = 2050
$ws.Range("L106").Value = 80027.7
$ws.Range("M106").Value = -1419
$ws.Range("N106").Value = -81289.7
$ws.Range("H123").Value = 37352
$ws.Range("J123").Value = 37352
$ws.Range("L123").Value = 37352
$ws.Range("N123").Value = -47152
$ws.Range("H128").Value = 46115
$ws.Range("J128").Value = 46115
$ws.Range("L128").Value = 46115
$ws.Range("N128").Value = -56075
$ws.Range("H130").Value = 47386
$ws.Range("J130").Value = 47386
$ws.Range("L130").Value = 47386
$ws.Range("N130").Value = -57426
$ws.Range("H133").Value = 50637.375
$ws.Range("J133").Value = 50637.375
$ws.Range("L133").Value = 50637.375
$ws.Range("N133").Value = -60757.375
$ws.Range("H135").Value = 15626285
$ws.Range("I135").Value = 870.0769
$ws.Range("J135").Value = 83336420
$ws.Range("K135").Value = 7830.6921
$ws.Range("L135").Value = 750027780
$ws.Range("M135").Value = -5295.6921
$ws.Range("N135").Value = -750032850
$ws.Range("H138").Value = 1517.6477
$ws.Range("I138").Value = 888.8431399999999
$ws.Range("J138").Value = 2384.3784
$ws.Range("K138").Value = 2666.52942
$ws.Range("L138").Value = 7153.135200000001
$ws.Range("M138").Value = 2473.47058
$ws.Range("N138").Value = -17433.1352

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10432.271
$ws.Range("I32").Value = 9749.98
$ws.Range("K32").Value = 9749.98
$ws.Range("M32").Value = -9462.98
$ws.Range("H103").Value = 40293.332
$ws.Range("J103").Value = 40293.332
$ws.Range("L103").Value = 40293.332
$ws.Range("N103").Value = -42637.332
$ws.Range("H128").Value = 50373
$ws.Range("J128").Value = 50373
$ws.Range("L128").Value = 50373
$ws.Range("N128").Value = -60333
$ws.Range("H131").Value = 44348.668
$ws.Range("J131").Value = 44348.668
$ws.Range("L131").Value = 44348.668
$ws.Range("N131").Value = -54428.668
$ws.Range("H132").Value = 3227.2666
$ws.Range("I132").Value = 1818.75
$ws.Range("K132").Value = 5456.25
$ws.Range("M132").Value = -2926.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1613.125
$ws.Range("I20").Value = 1137.5385
$ws.Range("J20").Value = 1938.5264
$ws.Range("K20").Value = 1137.5385
$ws.Range("L20").Value = 1938.5264
$ws.Range("M20").Value = -890.5385000000001
$ws.Range("N20").Value = -2432.5264
$ws.Range("H57").Value = 55369.5
$ws.Range("J57").Value = 55369.5
$ws.Range("L57").Value = 55369.5
$ws.Range("N57").Value = -56809.5
$ws.Range("H122").Value = 40334.8
$ws.Range("J122").Value = 40334.8
$ws.Range("L122").Value = 40334.8
$ws.Range("N122").Value = -50134.8
$ws.Range("H136").Value = 55369.5
$ws.Range("J136").Value = 55369.5
$ws.Range("L136").Value = 55369.5
$ws.Range("N136").Value = -65569.5
$ws.Range("H139").Value = 43935.8
$ws.Range("J139").Value = 43935.8
$ws.Range("L139").Value = 43935.8
$ws.Range("N139").Value = -54215.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H33").Value = 2240
$ws.Range("I33").Value = 1925
$ws.Range("J33").Value = 3500
$ws.Range("K33").Value = 1925
$ws.Range("L33").Value = 3500
$ws.Range("M33").Value = -1546
$ws.Range("N33").Value = -4258
$ws.Range("H52").Value = 69399.664
$ws.Range("J52").Value = 69399.664
$ws.Range("L52").Value = 69399.664
$ws.Range("N52").Value = -69987.664
$ws.Range("H105").Value = 3818.9092
$ws.Range("I105").Value = 3667.5557
$ws.Range("J105").Value = 4500
$ws.Range("K105").Value = 3667.5557
$ws.Range("L105").Value = 4500
$ws.Range("M105").Value = -1920.5557
$ws.Range("N105").Value = -7994
$ws.Range("H110").Value = 40424.8
$ws.Range("J110").Value = 40424.8
$ws.Range("L110").Value = 40424.8
$ws.Range("N110").Value = -48604.8
$ws.Range("H132").Value = 130417.82
$ws.Range("I132").Value = 2659.8
$ws.Range("J132").Value = 236882.83
$ws.Range("K132").Value = 7979.400000000001
$ws.Range("L132").Value = 710648.49
$ws.Range("M132").Value = -5449.400000000001
$ws.Range("N132").Value = -715708.49
$ws.Range("H137").Value = 36476.92
$ws.Range("J137").Value = 36476.92
$ws.Range("L137").Value = 36476.92
$ws.Range("N137").Value = -46676.92
$ws.Range("H139").Value = 60399.8
$ws.Range("J139").Value = 64999.75
$ws.Range("L139").Value = 64999.75
$ws.Range("N139").Value = -75279.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 3577.4856
$ws.Range("I113").Value = 6806.125
$ws.Range("J113").Value = 858.6316
$ws.Range("K113").Value = 20418.375
$ws.Range("L113").Value = 2575.8948
$ws.Range("M113").Value = -18248.375
$ws.Range("N113").Value = -6915.8948
$ws.Range("H124").Value = 1264.6296
$ws.Range("I124").Value = 976.6667
$ws.Range("J124").Value = 1300.625
$ws.Range("K124").Value = 2930.0001
$ws.Range("L124").Value = 3901.875
$ws.Range("M124").Value = 1979.9999
$ws.Range("N124").Value = -13721.875
$ws.Range("H134").Value = 33370198
$ws.Range("I134").Value = 43524204
$ws.Range("K134").Value = 130572612
$ws.Range("M134").Value = -130567542
$ws.Range("H140").Value = 168626.83
$ws.Range("I140").Value = 211919.89
$ws.Range("K140").Value = 635759.67
$ws.Range("M140").Value = -630579.67
$ws.Range("H141").Value = 83336230
$ws.Range("I141").Value = 111113630
$ws.Range("J141").Value = 4046.3333
$ws.Range("K141").Value = 333340890
$ws.Range("L141").Value = 12138.9999
$ws.Range("M141").Value = -333335710
$ws.Range("N141").Value = -22498.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H130").Value = 53992
$ws.Range("J130").Value = 53992
$ws.Range("L130").Value = 53992
$ws.Range("N130").Value = -64032
$ws.Range("H132").Value = 3091.8845
$ws.Range("I132").Value = 2226.4666
$ws.Range("J132").Value = 4272
$ws.Range("K132").Value = 6679.399800000001
$ws.Range("L132").Value = 12816
$ws.Range("M132").Value = -4149.399800000001
$ws.Range("N132").Value = -17876
$ws.Range("H137").Value = 41635.8
$ws.Range("J137").Value = 41635.8
$ws.Range("L137").Value = 41635.8
$ws.Range("N137").Value = -51835.8
$ws.Range("H139").Value = 33947
$ws.Range("J139").Value = 33947
$ws.Range("L139").Value = 33947
$ws.Range("N139").Value = -44227

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2975.3784
$ws.Range("I132").Value = 2100.85
$ws.Range("J132").Value = 4004.2354
$ws.Range("K132").Value = 6302.549999999999
$ws.Range("L132").Value = 12012.7062
$ws.Range("M132").Value = -3772.549999999999
$ws.Range("N132").Value = -17072.7062
$ws.Range("H133").Value = 39712.855
$ws.Range("J133").Value = 39712.855
$ws.Range("L133").Value = 39712.855
$ws.Range("N133").Value = -44772.855
$ws.Range("H134").Value = 50193.4
$ws.Range("J134").Value = 50193.4
$ws.Range("L134").Value = 50193.4
$ws.Range("N134").Value = -60333.4
$ws.Range("H137").Value = 41083.332
$ws.Range("J137").Value = 41083.332
$ws.Range("L137").Value = 41083.332
$ws.Range("N137").Value = -51283.332
$ws.Range("H139").Value = 50499.6
$ws.Range("J139").Value = 50499.6
$ws.Range("L139").Value = 50499.6
$ws.Range("N139").Value = -60779.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H125").Value = 41565.6
$ws.Range("J125").Value = 41565.6
$ws.Range("L125").Value = 41565.6
$ws.Range("N125").Value = -51405.6
$ws.Range("H132").Value = 1522.4667
$ws.Range("I132").Value = 1254.2188
$ws.Range("J132").Value = 2182.7693
$ws.Range("K132").Value = 3762.6564
$ws.Range("L132").Value = 6548.3079
$ws.Range("M132").Value = -1232.6564
$ws.Range("N132").Value = -11608.3079
$ws.Range("H139").Value = 51939.8
$ws.Range("J139").Value = 51939.8
$ws.Range("L139").Value = 51939.8
$ws.Range("N139").Value = -62219.8
